$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp banner in row 1 ---
$ws.Range("A1").Value2 = "Datos actualizados a 22 de Julio de 2020 a las 09:08"

# --- Refresh per-country stats that changed in this data pull ---
# Columns: B Casos totales, C Nuevos casos, D Casos activos, E Recuperados,
#          F Casos criticos, G Muertes hoy, H Muertes
function Set-CountryRow($country, $totalCases, $newCases, $activeCases, $recovered, $criticalCases, $deathsToday, $deaths) {
    $found = $ws.Range("A4:A219").Find($country)
    if ($found -eq $null) {
        Write-Output "WARNING: country not found: $country"
        return
    }
    $r = $found.Row
    $ws.Cells.Item($r, 2).Value2 = $totalCases    # B Casos totales
    $ws.Cells.Item($r, 3).Value2 = $newCases      # C Nuevos casos
    $ws.Cells.Item($r, 4).Value2 = $activeCases   # D Casos activos
    $ws.Cells.Item($r, 5).Value2 = $recovered     # E Recuperados
    $ws.Cells.Item($r, 6).Value2 = $criticalCases # F Casos criticos
    $ws.Cells.Item($r, 7).Value2 = $deathsToday   # G Muertes hoy
    $ws.Cells.Item($r, 8).Value2 = $deaths        # H Muertes
}

Set-CountryRow "Estados Unidos" 4028733 164 1886778 1996997 0 5  144958
Set-CountryRow "India"          1194888 803 753050  413067  0 1  28771
Set-CountryRow "Afganistan"     35727   112 23924   10613   0 4  1190
Set-CountryRow "Armenia"        35693   439 24766   10249   0 16 678
Set-CountryRow "El Salvador"    12582   0   6996    5223    0 11 363
Set-CountryRow "Hungria"        4366    19  3283    487     0 0  596
Set-CountryRow "Letonia"        1197    4   1045    121     0 0  31
Set-CountryRow "Georgia"        1073    24  907     150     0 0  16

# --- Re-rank the country table by total cases (column B), descending ---
$dataRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$dataRange.Sort($sortKey, 2)

# --- Break the tie between Groenlandia and Islas Malvinas (equal totals) ---
# so that Islas Malvinas now ranks just above Groenlandia.
$rowGroenlandia = $ws.Range("A4:A219").Find("Groenlandia").Row
$rowMalvinas = $ws.Range("A4:A219").Find("Islas Malvinas").Row

if ($rowGroenlandia -lt $rowMalvinas) {
    for ($c = 1; $c -le 8; $c++) {
        $tmp = $ws.Cells.Item($rowGroenlandia, $c).Value2
        $ws.Cells.Item($rowGroenlandia, $c).Value2 = $ws.Cells.Item($rowMalvinas, $c).Value2
        $ws.Cells.Item($rowMalvinas, $c).Value2 = $tmp
    }
}

Write-Output "Update complete"
